$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-29 down to 19-30
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44589
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100112028
$ws.Range("G18").Value = "Sandia"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 325
$ws.Range("L18").Value = 350
$ws.Range("M18").Value = 338
$ws.Range("N18").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 338
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
